$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - EU Flot
$ws.Range("C2").Value = 0.0972481420361959
$ws.Range("D2").Value = 0.09724811614937853

# Row 3 - EU Equity
$ws.Range("C3").Value = 0.09649764182944713
$ws.Range("D3").Value = 0.09649766910532863

# Row 4 - US Equity
$ws.Range("B4").Value = 0.42
$ws.Range("C4").Value = 0.1473185240024383
$ws.Range("D4").Value = 0.1473185371538656

# Row 5 - Greek Gov
$ws.Range("C5").Value = 0.09898494248164209
$ws.Range("D5").Value = 0.09898493358485361

# Row 6 - EU Corps
$ws.Range("C6").Value = 0.3482936652637426
$ws.Range("D6").Value = 0.3482936594876294

# Row 7 - EU Gov
$ws.Range("C7").Value = 0.09828707602528544
$ws.Range("D7").Value = 0.09828711179714611

# Row 8 - Cash
$ws.Range("C8").Value = 0.1133700083612487
$ws.Range("D8").Value = 0.1133699727217983
